$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.828.30"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.093.63"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.74"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0836"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "2.404.22"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "2.095.07"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "38.749.60"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("E28").Value = "  +8.59%  "
$ws.Range("E29").Value = "  +13.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("E34").Value = "  +4.47%  "
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "1.540.74"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0922"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("E46").Value = "  +8.25%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "2.290.63"
$ws.Range("E51").Value = "  +2.42%  "
